{"js": "// Wrap the \"{sequential_translation}: {id}\" subtitle paragraph with the\n// sequential_enabled config-value template tags:\n//   {#sequential_enabled} ... {/sequential_enabled}\n//\n// Target paragraph (style \"Subtitle\"/\"Untertitel\") currently reads:\n//   {sequential_translation}: {id}\n// and must become:\n//   {#sequential_enabled}{sequential_translation}: {id}{/sequential_enabled}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"text,style\");\n}\nawait context.sync();\n\nconst target = paragraphs.items.find(\n  (p) => p.style === \"Subtitle\" && p.text.indexOf(\"sequential_translation\") !== -1\n);\nif (!target) {\n  throw new Error(\"Could not find the sequential_translation subtitle paragraph\");\n}\n\n// Insert the opening tag as a brand-new run at the very start of the paragraph.\nconst openRun = target.insertText(\"{#sequential_enabled}\", \"Start\");\n// Force the new run to materialize an explicit (empty) run-properties element,\n// matching the sibling runs already in this paragraph (each of which carries\n// its own, currently-empty, <w:rPr/>).\nopenRun.font.bold = true;\nawait context.sync();\nopenRun.font.bold = false;\nawait context.sync();\n\n// Insert the closing tag as a brand-new run at the very end of the paragraph.\nconst closeRun = target.insertText(\"{/sequential_enabled}\", \"End\");\ncloseRun.font.bold = true;\nawait context.sync();\ncloseRun.font.bold = false;\nawait context.sync();\n", "ps1": "# Wrap the \"{sequential_translation}: {id}\" subtitle paragraph with the\n# sequential_enabled config-value template tags:\n#   {#sequential_enabled} ... {/sequential_enabled}\n#\n# Target paragraph (style \"Subtitle\") currently reads:\n#   {sequential_translation}: {id}\n# and must become:\n#   {#sequential_enabled}{sequential_translation}: {id}{/sequential_enabled}\n\n$d = $word.ActiveDocument\n\n# Locate the subtitle paragraph that carries the sequential id placeholders.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Style.NameLocal -eq \"Subtitle\" -and $p.Range.Text -like \"*sequential_translation*\") {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not find the sequential_translation subtitle paragraph\"\n}\n\n# Insert the opening tag as a brand-new run at the very start of the paragraph.\n$startR = $target.Range.Duplicate\n$startR.Collapse(1)  # wdCollapseStart\n$startR.InsertBefore(\"{#sequential_enabled}\")\n# Force the new run to materialize an explicit (empty) run-properties element,\n# matching the sibling runs already in this paragraph (each of which carries\n# its own, currently-empty, rPr).\n$startR.Font.Bold = 1\n$startR.Font.Bold = 0\n\n# Insert the closing tag as a brand-new run at the very end of the paragraph\n# (right before the paragraph mark).\n$endR = $target.Range\n$endR.Collapse(0)  # wdCollapseEnd\n$endR.InsertAfter(\"{/sequential_enabled}\")\n$endR.Font.Bold = 1\n$endR.Font.Bold = 0\n"}
